$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle/caratula shape (the one holding "Curso: I4051").
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame -and $cand.TextFrame.TextRange.Text.IndexOf("Curso: I4051") -ge 0) {
        $shp = $cand
        break
    }
}

$tr = $shp.TextFrame.TextRange

# 1) "Curso: I4051" -> "Curso: I4051 (Palazzo)"
$idx = $tr.Text.IndexOf("Curso: I4051")
$tr.Characters($idx + 1, 12).Text = "Curso: I4051 (Palazzo)"

# 2) "Elaborado por: Rodrigo Maranzana" paragraph -> reword its lead-in
#    to "Docente: ", leaving the trailing "Rodrigo Maranzana" run as
#    its own (untouched) run within the same paragraph.
$prefix = "Elaborado por: "
$idx2 = $tr.Text.IndexOf($prefix)
$tr.Characters($idx2 + 1, $prefix.Length).Text = "Docente: "

# 3) Drop the old standalone "Docente: Martín Palazzo" paragraph
#    (paragraph 4) and the blank paragraph right after it (paragraph
#    5); deleting from the back first keeps earlier indices stable.
#    The paragraph touched in step 2 (paragraph 3) then inherits
#    paragraph 4's endParaRPr automatically, matching the authored
#    edit.
$tr.Paragraphs(5, 1).Delete()
$tr.Paragraphs(4, 1).Delete()
